$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PMTestData")

$newRow = 16

$ws.Cells.Item($newRow, 1).Value = "test_configureParallel_ringing_to_IP"
$ws.Cells.Item($newRow, 2).Value = "number_initiate -number 70001..70003 -numbertype ex,extension -i -d 70001..70003 -l 1 --csp 0,ip_extension -i -d 70001..70003,70001-70003,70001,70002,70003,parallel_ringing -e -d 70001,ip_extension -e -d 70001..70003,extension -e -d 70001..70003,number_end -number 70001..70003 -numbertype ex"
$ws.Cells.Item($newRow, 3).Value = "Y"

$ws.Cells.Item($newRow, 2).WrapText = $true
$ws.Rows.Item($newRow).RowHeight = 72.5

$ws.Range("F$newRow").Select()
